# Hawaii overview workbook - convert numeric "count" cells that are meant to
# read as plain text (to match the rest of the factsheet, which already
# stores every other figure as literal text), and add the missing County
# "Total" row.
#
# Helper: write a value into a cell as literal TEXT even when the value looks
# like a number (e.g. "630"). Excel's normal type-inference would otherwise
# store it as a numeric <v> cell, so we flip the cell to the "Text" number
# format first, assign the string, then restore the cell's style to Normal
# so no stray number-format is left attached to the cell itself.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overall
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overall")
Set-TextValue $ws.Range("A2") "630"

# ---------------------------------------------------------------------
# County
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("County")
Set-TextValue $ws.Range("B2") "108"
Set-TextValue $ws.Range("B3") "365"
Set-TextValue $ws.Range("B4") "40"
Set-TextValue $ws.Range("B5") "117"

# Kalawao County row gets re-expressed as percentages / currency instead of
# bare zeros.
Set-TextValue $ws.Range("B6") "0.00%"
Set-TextValue $ws.Range("C6") "`$0"
Set-TextValue $ws.Range("D6") "0.00%"
Set-TextValue $ws.Range("E6") "0.00%"
Set-TextValue $ws.Range("F6") "0.00%"

# New "Total" row.
$ws.Range("A7").Value = "Total"
Set-TextValue $ws.Range("B7") "630"
Set-TextValue $ws.Range("C7") "`$1,081,113,515"
Set-TextValue $ws.Range("D7") "8.92%"
Set-TextValue $ws.Range("E7") "-17.37%"
Set-TextValue $ws.Range("F7") "70.00%"

# ---------------------------------------------------------------------
# Congressional District
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Congressional District")
Set-TextValue $ws.Range("B2") "298"
Set-TextValue $ws.Range("B3") "332"
Set-TextValue $ws.Range("B4") "630"

# ---------------------------------------------------------------------
# Size
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Size")
Set-TextValue $ws.Range("B2") "201"
Set-TextValue $ws.Range("B3") "175"
Set-TextValue $ws.Range("B4") "117"
Set-TextValue $ws.Range("B5") "42"
Set-TextValue $ws.Range("B6") "67"
Set-TextValue $ws.Range("B7") "28"
Set-TextValue $ws.Range("B8") "630"

# ---------------------------------------------------------------------
# Subsector
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Subsector")
Set-TextValue $ws.Range("B2") "78"
Set-TextValue $ws.Range("B3") "65"
Set-TextValue $ws.Range("B4") "72"
Set-TextValue $ws.Range("B5") "62"
Set-TextValue $ws.Range("B6") "7"
Set-TextValue $ws.Range("B7") "155"
Set-TextValue $ws.Range("B8") "5"
Set-TextValue $ws.Range("B9") "55"
Set-TextValue $ws.Range("B10") "3"
Set-TextValue $ws.Range("B11") "123"
Set-TextValue $ws.Range("B12") "5"
Set-TextValue $ws.Range("B13") "630"
